# ---------------------------------------------------------------------------
# Edit 1: fix the typo "Frontend Developmen" -> "Frontend Development"
# ---------------------------------------------------------------------------
$d = $word.ActiveDocument
$d.Content.Find.Execute("Frontend Developmen", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Frontend Development", 2)

# ---------------------------------------------------------------------------
# Edit 2: strip the leftover direct-formatting (font color #434343 and the
# explicit "no underline" marker) from eight bullet-list requirement items.
# The color/underline were stray direct formatting left over from a paste;
# removing them lets the paragraphs inherit normal (automatic) formatting.
# ---------------------------------------------------------------------------

function Strip-StrayFormatting($paragraphXml) {
    # Drop the explicit font-color override ...
    $paragraphXml = [regex]::Replace($paragraphXml, "<w:color[^/]*/>", "")
    # ... and the explicit "underline = none" override.
    $paragraphXml = [regex]::Replace($paragraphXml, "<w:u w:val=`"none`"/>", "")
    # Keep the xml:space="preserve" marker on the text run regardless of
    # whether WordOpenXML decided the round-tripped text still needs it.
    $paragraphXml = [regex]::Replace($paragraphXml, "<w:t>", "<w:t xml:space=`"preserve`">")
    # WordOpenXML drops the run's w:rsidR when exporting a bare Range; put it
    # back so the run keeps the same identity attributes it started with.
    $paragraphXml = [regex]::Replace($paragraphXml,
        "<w:r w:rsidDel=`"00000000`" w:rsidRPr=`"00000000`">",
        "<w:r w:rsidDel=`"00000000`" w:rsidR=`"00000000`" w:rsidRPr=`"00000000`">")
    return $paragraphXml
}

$targetParagraphs = @(
    "The feedback form should provide an auto-correct feature to help the user write grammatically correct feedback. ",
    "The system shall provide a way to provide feedback anonymously ",
    "The system should provide a way for the user provide feedback from a public, user profile or format provided by the feedback form",
    "The system should give the user a way to thread different types of feedback",
    "The system should return feedback to the sender in a timely manner once completed by the sending-party",
    "The system should provide a way to use previous feedback to help answer similar feedback ",
    "The system shall provide users a way to report accounts which don’t consistently meet the general guidelines listed via the user manual",
    "The system shall provide a way to penalize users who aren’t providing constructive feedback or are harassing students/professors/administrators. "
)

foreach ($p in $d.Paragraphs) {
    $paraText = $p.Range.Text.TrimEnd([char]13, [char]7)
    foreach ($target in $targetParagraphs) {
        if ($paraText -eq $target) {
            $rng = $p.Range
            $openXml = $rng.WordOpenXML
            if ($openXml -match "(<w:p [^>]*>.*?</w:p>)") {
                $cleanedParagraph = Strip-StrayFormatting $matches[1]
                $pkg = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
                       '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
                       '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
                       '<pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml">' +
                       $cleanedParagraph +
                       '</w:document></pkg:xmlData></pkg:part></pkg:package>'
                $rng.InsertXML($pkg)
            }
        }
    }
}
